# Switching from different sequences to one shared sequence: remove the
# SEQ_NAME column (G) from the "rights_and_functions" sheet entirely,
# including its header, its per-row sequence values (cds2db_in_seq /
# db_log_seq) and its header comment, shifting every column from H
# onward one place to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

# Remove the review comment attached to the SEQ_NAME header (G33) before
# the column shifts away from under it.
$comment = $ws.Range("G33").Comment
if ($comment -ne $null) {
    $comment.Delete() | Out-Null
}

# Delete the whole SEQ_NAME column; everything to its right (RIGHTS,
# GRANT_TARGET_USER, COPY_FUNC_SCRIPT_NAME, COPY_FUNC_NAME, SCHEMA_2,
# TABLE_POSTFIX_2, SCHEMA_3, TABLE_POSTFIX_3) shifts left by one column.
$ws.Columns("G").Delete() | Out-Null

# Leave the selection where the author last left it when saving.
$ws.Range("I47").Select() | Out-Null
